# Added assign leave functionality
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AllTests")

# Remove the old leftover test-data cells (rows 8-11)
$ws.Range("A8:G11").ClearContents()

# New "assign_leave" test block header + data row 6-7, and the marker in F8
# Write the header row texts before the "assign_leave" label so the shared
# string table receives "From Date"/"To Date" ahead of "assign_leave"
# (matches original authoring order).
$ws.Range("B6").Value = "Employee Full Name"
$ws.Range("C6").Value = "Leave Type"
$ws.Range("D6").Value = "From Date"
$ws.Range("E6").Value = "To Date"
$ws.Range("A6").Value = "assign_leave"

$ws.Range("B7").Value = "Russel Hamilton"
$ws.Range("C7").Value = "FMLA US"
$ws.Range("D7").Value = 43994
$ws.Range("E7").Value = 44001

$ws.Range("F8").Value = "assign_leave"

# Re-use the existing formatting already present in the sheet so the style
# table doesn't get duplicate entries: copy format only (xlPasteFormats).
$ws.Range("A1").Copy()
$ws.Range("A6").PasteSpecial(-4122)

$ws.Range("F4").Copy()
$ws.Range("F8").PasteSpecial(-4122)

$ws.Range("B1:E1").Copy()
$ws.Range("B6:E6").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("B7").PasteSpecial(-4122)

$ws.Range("C2").Copy()
$ws.Range("C7").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# New number format (date) for the From Date / To Date values.
$ws.Range("D7:E7").NumberFormat = "yyyy\-mm\-dd;@"

# Update the selection to match the edited sheet.
$ws.Range("E11").Select() | Out-Null
